# Update cryptocurrency price/volume/hour data per latest scrape
# (Updated symbol list on Sun Feb  5 18:00:44 UTC 2023 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.99%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "18"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.98%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "18"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.435"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.53%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "18"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08100"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.05%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "18"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.725"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.90%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "18"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.322"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.30%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "18"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.888"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.89%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "18"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.43%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "18"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9444"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.14%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "18"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1178"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-8.41%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "18"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1893"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.53%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "18"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09662"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.79%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "18"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04202"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "9.86%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "18"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1069"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.81%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "18"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001275"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.09%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "18"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006089"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.17%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "18"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.560"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.53%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "18"

# Row 19
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "18"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.759"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.17%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "18"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1364"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.19%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "18"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.19%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "18"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04391"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "18"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001241"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.47%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "18"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004308"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.52%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "18"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001239"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.50%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "18"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004018"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "31.95%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "18"

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "18"

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "18"

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "18"

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "18"

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "18"

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "18"

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "18"

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "18"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "18"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "18"

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "18"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02678"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.43%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "18"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05572"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.79%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "18"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007750"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.54%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "18"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009784"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.11%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "18"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1399"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.71%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "18"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002127"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.65%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "18"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009632"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.85%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "18"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007109"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.23%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "18"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000756"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.82%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "18"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003485"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.36%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "18"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002287"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.48%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "18"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002115"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.82%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "18"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002015"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.82%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "18"

